$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set transparency milestone marking for row 24 (Demonstrate proper transparent masking...)
$ws.Range("E24").Value = "III"
$ws.Range("F24").Value = "X"

# Set transparency milestone marking for row 68 (Transparent Blending...)
$ws.Range("E68").Value = "III"
$ws.Range("F68").Value = "X"

# Mark milestone II and III complete for rows 90-91
$ws.Range("D90").Value = "X"
$ws.Range("E90").Value = "X"
$ws.Range("D91").Value = "X"
$ws.Range("E91").Value = "X"

# Update view state (scroll position and selection)
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E91").Select()
